$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Default row height for the DATA sheet changes from 14.5 to 38.5 (customHeight)
$ws.StandardHeight = 38.5

# Row 3 was amazonTest / yes / firefox / 97.0 / Laptops -> becomes chrome / 98.0
$ws.Range("C3").Value = "chrome"
$ws.Range("D3").Value = "'98.0"

# Row 4 was amazonTest / yes / chrome / 98.0 / Laptops -> becomes edge / 98.0
$ws.Range("C4").Value = "edge"
$ws.Range("D4").Value = "'98.0"
